$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 206 (existing rows 206:224 shift down to 208:226).
$ws.Rows("206:207").Insert()

# New row 206: Damasco / Dina / Especial (San Felipe de Aconcagua)
$ws.Range("A206").Value = 6
$ws.Range("B206").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C206").Value = "Metropolitana"
$ws.Range("D206").Value = 44918
$ws.Range("D206").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E206").Value = 13
$ws.Range("F206").Value = "Fruta"
$ws.Range("G206").Value = 100103
$ws.Range("H206").Value = "Frutos de hueso (carozo)"
$ws.Range("I206").Value = 100103003
$ws.Range("J206").Value = "Damasco"
$ws.Range("K206").Value = "Dina"
$ws.Range("L206").Value = "Especial"
$ws.Range("M206").Value = 45
$ws.Range("N206").Value = 18000
$ws.Range("O206").Value = 18000
$ws.Range("P206").Value = 18000
$ws.Range("Q206").Value = "$/caja 16 kilos"
$ws.Range("R206").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S206").Value = 1125
$ws.Range("T206").Value = 16

# New row 207: Damasco / Dina / Primera (San Felipe de Aconcagua)
$ws.Range("A207").Value = 6
$ws.Range("B207").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C207").Value = "Metropolitana"
$ws.Range("D207").Value = 44918
$ws.Range("D207").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E207").Value = 13
$ws.Range("F207").Value = "Fruta"
$ws.Range("G207").Value = 100103
$ws.Range("H207").Value = "Frutos de hueso (carozo)"
$ws.Range("I207").Value = 100103003
$ws.Range("J207").Value = "Damasco"
$ws.Range("K207").Value = "Dina"
$ws.Range("L207").Value = "Primera"
$ws.Range("M207").Value = 60
$ws.Range("N207").Value = 15000
$ws.Range("O207").Value = 15000
$ws.Range("P207").Value = 15000
$ws.Range("Q207").Value = "$/caja 16 kilos"
$ws.Range("R207").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S207").Value = 938
$ws.Range("T207").Value = 16
